# Update export template: remove the "Actual Expended" columns for both the
# Personnel and Contract sections (v2024-08 export template refresh).
#
# Column AC ("Personnel_Actual_Expended" / "Actual Personnel Expenditures")
# and column AG ("Contract_Actual_Expended" / "Contract Actual Expended")
# are deleted outright (full column delete + shift-left), which pulls the
# FTE Count / Justification / Contract_Estimated_Expended / etc. columns
# left to fill the gap and drops the now-unused help text strings.
#
# Columns must be removed right-to-left so that deleting the first (AC)
# doesn't shift AG's position out from under the second delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("AG").Delete()
$ws.Columns("AC").Delete()

# Restore the reporter's working selection to cell B8 (first data-entry row).
[void]$ws.Range("B8").Select()
